$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.328.53'
$ws.Range("E2").Value = '  -4.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.759.94'
$ws.Range("E3").Value = '  -4.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.64'
$ws.Range("E6").Value = '  -2.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4275'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3606'
$ws.Range("E8").Value = '  -1.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07050'
$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8287'
$ws.Range("E10").Value = '  -4.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.08'
$ws.Range("E11").Value = '  -2.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.734.32'
$ws.Range("E12").Value = '  -5.82%  '

$ws.Range("E13").Value = '  -4.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.345'
$ws.Range("E14").Value = '  -2.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06820'
$ws.Range("E15").Value = '  -1.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.96'
$ws.Range("E17").Value = '  -2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008628'
$ws.Range("E18").Value = '  -3.23%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.89'
$ws.Range("E20").Value = '  -3.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.184.63'
$ws.Range("E21").Value = '  -5.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.976'
$ws.Range("E22").Value = '  -3.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.06'
$ws.Range("E23").Value = '  +2.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.963.92'
$ws.Range("E24").Value = '  -5.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.902'
$ws.Range("E25").Value = '  -4.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.62'
$ws.Range("E26").Value = '  -2.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.07'
$ws.Range("E27").Value = '  -4.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.016'
$ws.Range("E28").Value = '  -2.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.08'
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.664'
$ws.Range("E30").Value = '  -8.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08855'
$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7195'
$ws.Range("E32").Value = '  -3.94%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.295'
$ws.Range("E33").Value = '  -5.47%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.100'
$ws.Range("E34").Value = '  -2.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9999'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.687'
$ws.Range("E36").Value = '  -10.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.064'
$ws.Range("E37").Value = '  -2.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05085'
$ws.Range("E38").Value = '  -4.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01873'
$ws.Range("E39").Value = '  -3.32%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1598'
$ws.Range("E40").Value = '  -3.33%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4879'
$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.476'
$ws.Range("E42").Value = '  -11.56%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.129'
$ws.Range("E43").Value = '  -5.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.972'
$ws.Range("E44").Value = '  -4.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.17'
$ws.Range("E45").Value = '  -1.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.02'
$ws.Range("E47").Value = '  -4.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06172'
$ws.Range("E48").Value = '  -4.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4456'
$ws.Range("E49").Value = '  -4.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.563'
$ws.Range("E50").Value = '  -3.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.703'
$ws.Range("E51").Value = '  -1.30%  '
